# Add a new "Save" column (H) to the s_vals sheet, matching the style of
# the existing header row and populating the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: same bold/border/centered style as the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data cells H2:H3 (plain numeric, no special style).
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
